$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column Q (30-jun) with its hourly prices.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Clone the header formatting from the previous day's column (P1) onto the
# new header cell (Q1) before writing its value.
$wsPrix.Range("P1").Copy()
$wsPrix.Range("Q1").PasteSpecial(-4122) # xlPasteFormats
$wsPrix.Range("Q1").Value = "30-jun"

$prixValues = @(
    98.73,
    89.03,
    85.03,
    84.13,
    82.98,
    82.02,
    97.73,
    111.68,
    108.58,
    92.42,
    79.83,
    69.05,
    52.44,
    46.31,
    51.53,
    64.76000000000001,
    83.83,
    93.56999999999999,
    110.65,
    178.94,
    185,
    175.01,
    157,
    122.65
)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 17).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append rows for 2025-06-28 and 2025-06-29.
# The Date column holds plain text (e.g. "2025-06-28"), not real Excel
# dates, so force a text number format while writing the value, then drop
# back to the default "Normal" style to avoid leaving any formatting on
# the new cells (matching the unstyled cells used by the existing rows).
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A14").NumberFormat = "@"
$wsGaz.Range("A14").Value = "2025-06-28"
$wsGaz.Range("A14").Style = "Normal"
$wsGaz.Range("B14").Value = 32.675

$wsGaz.Range("A15").NumberFormat = "@"
$wsGaz.Range("A15").Value = "2025-06-29"
$wsGaz.Range("A15").Style = "Normal"
$wsGaz.Range("B15").Value = 32.675

# ---------------------------------------------------------------------------
# Sheet "CO2": append rows for 2025-06-28 and 2025-06-29.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A14").NumberFormat = "@"
$wsCo2.Range("A14").Value = "2025-06-28"
$wsCo2.Range("A14").Style = "Normal"
$wsCo2.Range("B14").Value = 69.92

$wsCo2.Range("A15").NumberFormat = "@"
$wsCo2.Range("A15").Value = "2025-06-29"
$wsCo2.Range("A15").Style = "Normal"
$wsCo2.Range("B15").Value = 69.92
